$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 updates
$ws.Range("B5").Value = "Desh mobile"
$ws.Range("C5").Value = "c25s-128"
$ws.Range("D5").Value = 867623053838711
$ws.Range("F5").Value = "30.11.2021"

# Row 6 updates
$ws.Range("B6").Value = "Shohan enterprise "
$ws.Range("D6").Value = 869092052729339
$ws.Range("F6").Value = "30.11.2021"

# Row 7 updates
$ws.Range("B7").Value = "Apurbo mobile"
$ws.Range("C7").Value = "c20a"
$ws.Range("D7").Value = 868790052445811
$ws.Range("F7").Value = "30.11.2021"

# Row 8 - new data
$ws.Range("B8").Value = "Apurbo mobile"
$ws.Range("C8").Value = "C21/64"
$ws.Range("D8").Value = 864623050833338
$ws.Range("E8").Value = "P"
$ws.Range("F8").Value = "30.11.2021"

# Row 9 - new data
$ws.Range("B9").Value = "friends mobile"
$ws.Range("D9").Value = 867623050430934
$ws.Range("E9").Value = "P+C"
$ws.Range("C9").Value = "C25s/64"
$ws.Range("F9").Value = "30.11.2021"

# Selection change noted in diff
$ws.Range("I8").Select()
